# Adds a new "DC Unit Loading Details" column (E1:E3) with header + two
# labels to both loop sheets, copying the existing header/label cell
# formatting (E1 like the table header row, E2/E3 like the table label
# cells) and updates the selection on each sheet to E1:E3.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Add_Devices_LoopA", "Add_Devices_LoopB")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # New labels in column E, rows 1-3.
    $ws.Range("E1").Value = "DC Unit Loading Details Name"
    $ws.Range("E2").Value = "Current (DC Units)"
    $ws.Range("E3").Value = "Current (worst case)"

    # Match formatting: E1 picks up the bold header style (same as A7:J7),
    # E2/E3 pick up the shaded label style (same as A8:J8).
    $ws.Range("A7").Copy()
    $ws.Range("E1").PasteSpecial(-4122)

    $ws.Range("A8").Copy()
    $ws.Range("E2:E3").PasteSpecial(-4122)

    $ws.Range("E1:E3").Select()
}

# Restore the originally active sheet/tab (Add_Devices_LoopB) and its
# selection, since selecting on LoopA above switches the active sheet.
$wsB = $wb.Worksheets.Item("Add_Devices_LoopB")
$wsB.Activate()
$wsB.Range("E1:E3").Select()

$excel.CutCopyMode = $false
